$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values - use apostrophe prefix to force text,
# then reset style to Normal so no residual NumberFormat/style is left on the cell,
# preserving the original inlineStr-as-text representation.
$ws.Range("D2").Value = "'66.388.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.444.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'599.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'146.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.444.91"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.135"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'6.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "'0.0000213"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'4.021.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'30.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.445.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'66.321.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'14.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'436.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'8.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.616"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'76.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'3.578.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0000122"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Value = "'8.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Value = "'0.159"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'25.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Value = "'1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'3.425.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'7.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'172.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'0.0852"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'5.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.870"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'45.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'25.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'7.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'2.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.988"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (E) column values - plain text, never numeric-looking so no
# special handling required.
$ws.Range("E2").Value = "  -5.00%  "
$ws.Range("E3").Value = "  -6.67%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -7.49%  "
$ws.Range("E6").Value = "  -9.39%  "
$ws.Range("E7").Value = "  -6.61%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -5.35%  "
$ws.Range("E10").Value = "  -7.46%  "
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("E12").Value = "  -6.42%  "
$ws.Range("E13").Value = "  -8.72%  "
$ws.Range("E14").Value = "  -6.81%  "
$ws.Range("E15").Value = "  -6.13%  "
$ws.Range("E16").Value = "  -5.72%  "
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  -7.90%  "
$ws.Range("E21").Value = "  -7.49%  "
$ws.Range("E22").Value = "  -14.42%  "
$ws.Range("E23").Value = "  -5.68%  "
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -6.78%  "
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("E28").Value = "  -9.20%  "
$ws.Range("E29").Value = "  -11.02%  "
$ws.Range("E30").Value = "  -6.57%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -11.96%  "
$ws.Range("E33").Value = "  -5.93%  "
$ws.Range("E34").Value = "  -5.45%  "
$ws.Range("E35").Value = "  -7.71%  "
$ws.Range("E36").Value = "  -9.09%  "
$ws.Range("E37").Value = "  -7.14%  "
$ws.Range("E38").Value = "  -7.71%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("E42").Value = "  -6.10%  "
$ws.Range("E43").Value = "  -5.76%  "
$ws.Range("E44").Value = "  -9.40%  "
$ws.Range("E45").Value = "  -6.76%  "
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("E48").Value = "  -12.61%  "
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("E50").Value = "  -14.90%  "
$ws.Range("E51").Value = "  -6.48%  "
